$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that used to sit right
#    under the H1 title at the top of the document.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play Fetching Fruits for Free - Review"
#    right before the final ("Prompt for DALLE ...") paragraph.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($n)
$insertStart = $pLast.Range.Start
$insertPoint = $d.Range($insertStart, $insertStart)

$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fetching Fruits for Free - Review</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xml) | Out-Null

# InsertXML leaves one spare empty paragraph behind, between the newly
# inserted paragraph and the "Prompt for DALLE ..." paragraph - drop it.
$spare = $d.Paragraphs.Item($n + 1)
$spare.Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 3) Swap out the DALLE prompt text for the meta-description copy, keeping
#    the run's existing (italic) formatting.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Prompt for DALLE: Create a feature image for Fetching Fruits game that follows the given criteria: - The image should be in cartoon style. - The image should feature a happy Maya warrior with glasses. The image should be bright and vibrant, featuring the happy Maya warrior with glasses in the center of the frame. The warrior should be holding a basket of colorful fruits, with a big smile on his face. He should be surrounded by various fruit symbols from the game, such as cherries, lemons, oranges, plums, watermelons, and bells. The background should have a jungle feel to it, with lush green trees and vines. The overall style of the image should be fun and energetic, conveying the excitement of the game.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Fetching Fruits, a free slot game with smooth gameplay mechanics, flexible betting options, and an autoplay feature.", 2) | Out-Null
